# The test fixture workbook gets a block of 15 rows x 9 columns of
# sample numbers (1..9 repeated across columns A:I) inserted above the
# existing chart source-data table, the leftover "A2" label (and its
# now-unused shared string) is removed, and the chart is shifted down
# to make room for the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old standalone "A2" label in A2 - it is no longer needed.
$ws.Range("A2").ClearContents() | Out-Null

# Fill rows 1-15, columns A-I (1-9) with the sequence 1..9 on every row.
for ($r = 1; $r -le 15; $r++) {
    for ($c = 1; $c -le 9; $c++) {
        $ws.Cells.Item($r, $c).Value = $c
    }
}

# Move the chart down so it still starts right below the new data block
# (it keeps its original size, just shifts down by roughly one row).
$co = $ws.ChartObjects().Item(1)
$co.Left = 12.75
$co.Top = 24.37496062992126
$co.Width = 386.625
$co.Height = 216

# Update the active selection like in the saved workbook.
$ws.Range("K11").Select() | Out-Null
